$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Price column (D) ---
# Some of the new values look like plain numbers (e.g. "0.4974"), and the
# COM layer will silently convert such strings into numeric cells when
# assigned directly. To keep them as text (matching the original
# inline/shared-string text cells) we force a Text number format on each
# target cell first, write the value, then reset the style back to
# "Normal" so no stray style attribute is left on the cell. NumberFormat
# and Style have to be applied to each cell individually (a joined/union
# range address does not get the formatting applied reliably).
$dCells = @("D2","D3","D5","D7","D8","D9","D11","D12","D14","D15","D16","D18","D19","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D44","D45","D46","D47","D48","D49","D51")
foreach ($cell in $dCells) {
    $ws.Range($cell).NumberFormat = "@"
}

$ws.Range("D2").Value = '27.922.00'
$ws.Range("D3").Value = '1.813.17'
$ws.Range("D5").Value = '311.01'
$ws.Range("D7").Value = '0.4974'
$ws.Range("D8").Value = '0.3953'
$ws.Range("D9").Value = '0.09593'
$ws.Range("D11").Value = '40.88'
$ws.Range("D12").Value = '6.426'
$ws.Range("D14").Value = '20.46'
$ws.Range("D15").Value = '1.816.56'
$ws.Range("D16").Value = '7.285'
$ws.Range("D18").Value = '92.38'
$ws.Range("D19").Value = '0.06657'
$ws.Range("D21").Value = '17.13'
$ws.Range("D22").Value = '5.910'
$ws.Range("D23").Value = '27.979.82'
$ws.Range("D25").Value = '2.252'
$ws.Range("D26").Value = '159.35'
$ws.Range("D27").Value = '2.021.08'
$ws.Range("D28").Value = '20.56'
$ws.Range("D29").Value = '2.383'
$ws.Range("D30").Value = '127.97'
$ws.Range("D31").Value = '0.1067'
$ws.Range("D32").Value = '1.035'
$ws.Range("D33").Value = '5.563'
$ws.Range("D34").Value = '3.636'
$ws.Range("D35").Value = '0.06708'
$ws.Range("D36").Value = '8.938'
$ws.Range("D37").Value = '0.02329'
$ws.Range("D38").Value = '0.2138'
$ws.Range("D39").Value = '4.938'
$ws.Range("D40").Value = '11.21'
$ws.Range("D44").Value = '13.10'
$ws.Range("D45").Value = '1.291'
$ws.Range("D46").Value = '0.5883'
$ws.Range("D47").Value = '3.696'
$ws.Range("D48").Value = '123.12'
$ws.Range("D49").Value = '1.933'
$ws.Range("D51").Value = '0.06777'

foreach ($cell in $dCells) {
    $ws.Range($cell).Style = "Normal"
}

# --- Update the Volume(1h) column (E) ---
# These values contain leading/trailing spaces and a percent sign, so
# Excel always keeps them as plain text; no special handling required.
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("E3").Value = '  +1.79%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.14%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("E7").Value = '  -3.18%  '
$ws.Range("E8").Value = '  +4.28%  '
$ws.Range("E9").Value = '  +23.24%  '
$ws.Range("E10").Value = '  +1.60%  '
$ws.Range("E11").Value = '  -1.07%  '
$ws.Range("E12").Value = '  +3.56%  '
$ws.Range("E13").Value = '  +0.08%  '
$ws.Range("E15").Value = '  +2.00%  '
$ws.Range("E16").Value = '  +1.71%  '
$ws.Range("E17").Value = '  +4.74%  '
$ws.Range("E18").Value = '  +0.98%  '
$ws.Range("E19").Value = '  +1.53%  '
$ws.Range("E20").Value = '  -0.02%  '
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("E22").Value = '  -0.04%  '
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("E24").Value = '  +1.36%  '
$ws.Range("E25").Value = '  +0.70%  '
$ws.Range("E26").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +1.67%  '
$ws.Range("E29").Value = '  +1.02%  '
$ws.Range("E30").Value = '  +2.37%  '
$ws.Range("E31").Value = '  -0.49%  '
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("E34").Value = '  +0.10%  '
$ws.Range("E35").Value = '  -5.05%  '
$ws.Range("E36").Value = '  +2.73%  '
$ws.Range("E37").Value = '  +0.56%  '
$ws.Range("E38").Value = '  +0.75%  '
$ws.Range("E39").Value = '  -1.25%  '
$ws.Range("E40").Value = '  -2.71%  '
$ws.Range("E41").Value = '  +1.41%  '
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("E44").Value = '  -0.19%  '
$ws.Range("E45").Value = '  -2.21%  '
$ws.Range("E46").Value = '  -0.61%  '
$ws.Range("E47").Value = '  -0.38%  '
$ws.Range("E48").Value = '  -3.70%  '
$ws.Range("E49").Value = '  +2.11%  '
$ws.Range("E50").Value = '  -2.26%  '
$ws.Range("E51").Value = '  -0.15%  '
